$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "62.838.42"
$ws.Range("E2").Value = "  +6.05%  "
$ws.Range("D3").Value = "2.445.60"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  +0.20%  "
Set-TextValue "D5" "580.18"
$ws.Range("E5").Value = "  +4.34%  "
Set-TextValue "D6" "145.68"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("D9").Value = "2.443.59"
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("E10").Value = "  +6.10%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("E13").Value = "  +5.14%  "
$ws.Range("E14").Value = "  +6.37%  "
$ws.Range("E15").Value = "  +9.40%  "
$ws.Range("D16").Value = "2.896.58"
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("D17").Value = "62.701.10"
$ws.Range("E17").Value = "  +5.89%  "
$ws.Range("D18").Value = "2.456.28"
$ws.Range("E18").Value = "  +4.61%  "
Set-TextValue "D19" "7.96"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +5.24%  "
Set-TextValue "D21" "326.51"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("E23").Value = "  +14.56%  "
$ws.Range("E24").Value = "  -0.10%  "
Set-TextValue "D25" "65.74"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("B26").Value = "Bittensor"
$ws.Range("C26").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D26" "616.00"
$ws.Range("E26").Value = "  +12.11%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D27" "1.10"
$ws.Range("E27").Value = "  +10.24%  "
$ws.Range("E28").Value = "  +4.23%  "
$ws.Range("D29").Value = "0.0₃0981"
$ws.Range("E29").Value = "  +8.33%  "
$ws.Range("D30").Value = "2.565.05"
Set-TextValue "D31" "8.13"
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("E32").Value = "  +9.72%  "
Set-TextValue "D33" "1.84"
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("E34").Value = "  +6.02%  "
$ws.Range("B35").Value = "BabyDogeCoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D35").Value = "0.0₆0379"
$ws.Range("E35").Value = "  +34.29%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "1.48"
$ws.Range("E36").Value = "  +5.64%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D38" "4.75"
$ws.Range("E38").Value = "  +5.59%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D39" "0.372"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D40" "152.26"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D41" "5.38"
$ws.Range("E41").Value = "  +8.05%  "
$ws.Range("B42").Value = "EthereumClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D42" "18.59"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.74"
$ws.Range("E43").Value = "  +18.15%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.76"
$ws.Range("E44").Value = "  +8.32%  "
Set-TextValue "D45" "42.37"
$ws.Range("E45").Value = "  +3.14%  "
Set-TextValue "D47" "144.07"
$ws.Range("E47").Value = "  +4.14%  "
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("E49").Value = "  +6.48%  "
$ws.Range("E50").Value = "  +2.98%  "
$ws.Range("E51").Value = "  +3.58%  "
